{"js": "// Applies the diff: removes stray <w:pPr>/<w:rFonts hint=\"eastAsia\"/> wrappers\n// from several paragraphs (turning some into bare empty <w:p/> elements),\n// strips the <w:rPr> from the \"commoniodemo\" list paragraph's <w:pPr>, and\n// inserts a new \"springmybatisdemo\" block (5 paragraphs) near the end of the\n// document, right before the final bookmark (\"_GoBack\") paragraph.\n\nconst OOXML_NS =\n  '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>{0}</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nfunction pkg(bodyInnerXml) {\n  return OOXML_NS.replace(\"{0}\", bodyInnerXml);\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraphs by their (unique) text content rather than by\n// a hard-coded index, so the script is resilient to minor structural drift.\nlet idxC3p3 = -1;\nlet idxCommonIoDemo = -1;\nlet idxAfterCommonsIoBlankPara = -1;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"c3p3\u6570\u636e\u5e93\u8fde\u63a5\u6c60\u5355\u72ec\u4f7f\u7528\u793a\u4f8b\") idxC3p3 = i;\n  if (t === \"commoniodemo\") idxCommonIoDemo = i;\n}\n\nif (idxC3p3 === -1) throw new Error(\"Could not locate 'c3p3' paragraph\");\nif (idxCommonIoDemo === -1) throw new Error(\"Could not locate 'commoniodemo' paragraph\");\n\n// The two blank paragraphs right after the c3p3 paragraph (each currently\n// only carrying a <w:pPr><w:rPr><w:rFonts hint=\"eastAsia\"/></w:rPr></w:pPr>).\nconst idxBlank1 = idxC3p3 + 1;\nconst idxBlank2 = idxC3p3 + 2;\n\n// Sanity-check: blank1/blank2 should indeed sit right before idxCommonIoDemo.\nif (idxBlank2 + 1 !== idxCommonIoDemo) {\n  throw new Error(\"Unexpected document structure around c3p3/commoniodemo paragraphs\");\n}\n\n// \"\u4f7f\u7528commonsio\u8bfb\u53d6\u6587\u4ef6\u4e0e\u5199\u5165\u6587\u4ef6\u793a\u4f8b\" is right after commoniodemo, and the\n// trailing bare <w:p/> is right after that. The new block is inserted right\n// after that bare paragraph (i.e. right before the final bookmark paragraph).\nconst idxCommonsioText = idxCommonIoDemo + 1;\nidxAfterCommonsIoBlankPara = idxCommonsioText + 1;\n\nif (paragraphs.items[idxCommonsioText].text !== \"\u4f7f\u7528commonsio\u8bfb\u53d6\u6587\u4ef6\u4e0e\u5199\u5165\u6587\u4ef6\u793a\u4f8b\") {\n  throw new Error(\"Could not locate '\u4f7f\u7528commonsio\u8bfb\u53d6\u6587\u4ef6\u4e0e\u5199\u5165\u6587\u4ef6\u793a\u4f8b' paragraph\");\n}\nif (paragraphs.items[idxAfterCommonsIoBlankPara].text !== \"\") {\n  throw new Error(\"Could not locate trailing blank paragraph after commonsio text\");\n}\n\n// --- Operation 1: strip <w:pPr> from the \"c3p3 / \u6570\u636e\u5e93\u8fde\u63a5\u6c60\u5355\u72ec\u4f7f\u7528\u793a\u4f8b\" paragraph ---\n{\n  const xml = pkg(\n    \"<w:p>\" +\n      '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>c3p3</w:t></w:r>' +\n      '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\\u6570\\u636e\\u5e93\\u8fde\\u63a5\\u6c60\\u5355\\u72ec\\u4f7f\\u7528\\u793a\\u4f8b</w:t></w:r>' +\n      \"</w:p>\"\n  );\n  paragraphs.items[idxC3p3].getRange().insertOoxml(xml, Word.InsertLocation.replace);\n}\n\n// --- Operation 2 & 3: collapse the two following blank paragraphs to bare <w:p/> ---\n{\n  const xml = pkg(\"<w:p/>\");\n  paragraphs.items[idxBlank1].getRange().insertOoxml(xml, Word.InsertLocation.replace);\n}\n{\n  const xml = pkg(\"<w:p/>\");\n  paragraphs.items[idxBlank2].getRange().insertOoxml(xml, Word.InsertLocation.replace);\n}\n\n// --- Operation 4: drop the <w:rPr> from \"commoniodemo\"'s <w:pPr> (keep pStyle/numPr) ---\n{\n  const xml = pkg(\n    \"<w:p>\" +\n      \"<w:pPr>\" +\n      '<w:pStyle w:val=\"2\"/>' +\n      '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n      \"</w:pPr>\" +\n      \"<w:r><w:t>commoniodemo</w:t></w:r>\" +\n      \"</w:p>\"\n  );\n  paragraphs.items[idxCommonIoDemo].getRange().insertOoxml(xml, Word.InsertLocation.replace);\n}\n\n// --- Operation 5: insert the new \"springmybatisdemo\" block after the trailing\n// blank paragraph that follows \"\u4f7f\u7528commonsio\u8bfb\u53d6\u6587\u4ef6\u4e0e\u5199\u5165\u6587\u4ef6\u793a\u4f8b\" ---\n{\n  const rFontsRPr = '<w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr>';\n  const runTexts = [\n    \"spring\",\n    \"\\u7ed3\\u5408\",\n    \"springmvc\",\n    \"\\u4e0e\",\n    \"mybatis\",\n    \"\\u7ed3\\u5408\\u6ca1\\u6709\\u5206\\u9875\\u5982\\u679c\\u4f7f\\u7528\\u5206\\u9875\\u8bf7\\u53c2\\u7167\",\n    \"mybatis\",\n    \"\\u5206\\u9875\\u793a\\u4f8b\",\n  ];\n  const runs = runTexts\n    .map((t) => `<w:r>${rFontsRPr}<w:t>${t}</w:t></w:r>`)\n    .join(\"\");\n\n  const newBlock =\n    // a. empty paragraph\n    `<w:p><w:pPr>${rFontsRPr}</w:pPr></w:p>` +\n    // b. \"springmybatisdemo\" list paragraph\n    \"<w:p>\" +\n    \"<w:pPr>\" +\n    '<w:pStyle w:val=\"2\"/>' +\n    '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n    rFontsRPr +\n    \"</w:pPr>\" +\n    \"<w:r><w:t>springmybatisdemo</w:t></w:r>\" +\n    \"</w:p>\" +\n    // c. description paragraph with the multiple runs\n    `<w:p><w:pPr>${rFontsRPr}</w:pPr>${runs}</w:p>` +\n    // d. empty paragraph\n    `<w:p><w:pPr>${rFontsRPr}</w:pPr></w:p>` +\n    // e. empty paragraph\n    `<w:p><w:pPr>${rFontsRPr}</w:pPr></w:p>`;\n\n  const xml = pkg(newBlock);\n  paragraphs.items[idxAfterCommonsIoBlankPara]\n    .getRange()\n    .insertOoxml(xml, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Applies the diff: removes stray <w:pPr>/<w:rFonts hint=\"eastAsia\"/> wrappers\n# from several paragraphs (turning some into bare empty <w:p/> elements),\n# strips the <w:rPr> from the \"commoniodemo\" list paragraph's <w:pPr>, and\n# inserts a new \"springmybatisdemo\" block (5 paragraphs) near the end of the\n# document, right before the final bookmark (\"_GoBack\") paragraph.\n\n$d = $word.ActiveDocument\n\nfunction CleanText($t) {\n    return $t.TrimEnd([char]13, [char]7)\n}\n\n# Locate the anchor paragraphs by their (unique) text content rather than a\n# hard-coded index, so the script is resilient to minor structural drift.\n$idxC3p3 = -1\n$idxCommonIoDemo = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = CleanText($d.Paragraphs.Item($i).Range.Text)\n    if ($t -eq \"c3p3\u6570\u636e\u5e93\u8fde\u63a5\u6c60\u5355\u72ec\u4f7f\u7528\u793a\u4f8b\") { $idxC3p3 = $i }\n    if ($t -eq \"commoniodemo\") { $idxCommonIoDemo = $i }\n}\n\nif ($idxC3p3 -eq -1) { throw \"Could not locate 'c3p3' paragraph\" }\nif ($idxCommonIoDemo -eq -1) { throw \"Could not locate 'commoniodemo' paragraph\" }\n\n# The two blank paragraphs right after the c3p3 paragraph (each currently\n# only carrying a <w:pPr><w:rPr><w:rFonts hint=\"eastAsia\"/></w:rPr></w:pPr>).\n$idxBlank1 = $idxC3p3 + 1\n$idxBlank2 = $idxC3p3 + 2\n\nif (($idxBlank2 + 1) -ne $idxCommonIoDemo) {\n    throw \"Unexpected document structure around c3p3/commoniodemo paragraphs\"\n}\n\n# \"\u4f7f\u7528commonsio\u8bfb\u53d6\u6587\u4ef6\u4e0e\u5199\u5165\u6587\u4ef6\u793a\u4f8b\" is right after commoniodemo, and the\n# trailing bare <w:p/> is right after that. The new block is inserted right\n# after that bare paragraph (i.e. right before the final bookmark paragraph).\n$idxCommonsioText = $idxCommonIoDemo + 1\n$idxAfterCommonsIoBlankPara = $idxCommonsioText + 1\n\nif ((CleanText($d.Paragraphs.Item($idxCommonsioText).Range.Text)) -ne \"\u4f7f\u7528commonsio\u8bfb\u53d6\u6587\u4ef6\u4e0e\u5199\u5165\u6587\u4ef6\u793a\u4f8b\") {\n    throw \"Could not locate '\u4f7f\u7528commonsio\u8bfb\u53d6\u6587\u4ef6\u4e0e\u5199\u5165\u6587\u4ef6\u793a\u4f8b' paragraph\"\n}\nif ((CleanText($d.Paragraphs.Item($idxAfterCommonsIoBlankPara).Range.Text)) -ne \"\") {\n    throw \"Could not locate trailing blank paragraph after commonsio text\"\n}\n\n$pkgHeader = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# --- Operation 1: strip <w:pPr> from the \"c3p3 / \u6570\u636e\u5e93\u8fde\u63a5\u6c60\u5355\u72ec\u4f7f\u7528\u793a\u4f8b\" paragraph ---\n$xml1 = $pkgHeader + `\n    '<w:p>' + `\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>c3p3</w:t></w:r>' + `\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\u6570\u636e\u5e93\u8fde\u63a5\u6c60\u5355\u72ec\u4f7f\u7528\u793a\u4f8b</w:t></w:r>' + `\n    '</w:p>' + `\n    $pkgFooter\n$d.Paragraphs.Item($idxC3p3).Range.InsertXML($xml1, \"Replace\") | Out-Null\n\n# --- Operation 2 & 3: collapse the two following blank paragraphs to bare <w:p/> ---\n$xmlBlank = $pkgHeader + '<w:p/>' + $pkgFooter\n$d.Paragraphs.Item($idxBlank1).Range.InsertXML($xmlBlank, \"Replace\") | Out-Null\n$d.Paragraphs.Item($idxBlank2).Range.InsertXML($xmlBlank, \"Replace\") | Out-Null\n\n# --- Operation 4: drop the <w:rPr> from \"commoniodemo\"'s <w:pPr> (keep pStyle/numPr) ---\n$xml4 = $pkgHeader + `\n    '<w:p>' + `\n    '<w:pPr>' + `\n    '<w:pStyle w:val=\"2\"/>' + `\n    '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' + `\n    '</w:pPr>' + `\n    '<w:r><w:t>commoniodemo</w:t></w:r>' + `\n    '</w:p>' + `\n    $pkgFooter\n$d.Paragraphs.Item($idxCommonIoDemo).Range.InsertXML($xml4, \"Replace\") | Out-Null\n\n# --- Operation 5: insert the new \"springmybatisdemo\" block after the trailing\n# blank paragraph that follows \"\u4f7f\u7528commonsio\u8bfb\u53d6\u6587\u4ef6\u4e0e\u5199\u5165\u6587\u4ef6\u793a\u4f8b\" ---\n$rFontsRPr = '<w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr>'\n$runTexts = @(\"spring\", \"\u7ed3\u5408\", \"springmvc\", \"\u4e0e\", \"mybatis\", \"\u7ed3\u5408\u6ca1\u6709\u5206\u9875\u5982\u679c\u4f7f\u7528\u5206\u9875\u8bf7\u53c2\u7167\", \"mybatis\", \"\u5206\u9875\u793a\u4f8b\")\n$runs = \"\"\nforeach ($t in $runTexts) {\n    $runs += \"<w:r>$rFontsRPr<w:t>$t</w:t></w:r>\"\n}\n\n$newBlock = `\n    \"<w:p><w:pPr>$rFontsRPr</w:pPr></w:p>\" + `\n    \"<w:p>\" + `\n    \"<w:pPr>\" + `\n    '<w:pStyle w:val=\"2\"/>' + `\n    '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' + `\n    $rFontsRPr + `\n    \"</w:pPr>\" + `\n    \"<w:r><w:t>springmybatisdemo</w:t></w:r>\" + `\n    \"</w:p>\" + `\n    \"<w:p><w:pPr>$rFontsRPr</w:pPr>$runs</w:p>\" + `\n    \"<w:p><w:pPr>$rFontsRPr</w:pPr></w:p>\" + `\n    \"<w:p><w:pPr>$rFontsRPr</w:pPr></w:p>\"\n\n$xml5 = $pkgHeader + $newBlock + $pkgFooter\n$d.Paragraphs.Item($idxAfterCommonsIoBlankPara).Range.InsertXML($xml5, \"After\") | Out-Null\n"}
